$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Extend the "System Design" paragraph: "...work together." ->
#    "...work together and reduce fallback." split across several runs,
#    mirroring how Word breaks runs while the text is being typed.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" It is important to understand the full system requirements in order to make it easier for members of the team to effectively work together.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$rng.Text = " It is important to understand the full system requirements in order to make it easier for members of the te"
$rng.Collapse(0)
$rng.InsertAfter("am to effectively work together and reduce ")
$rng.Font.NameBi = "Arial"

$rng.Collapse(0)
$rng.InsertAfter("fallback")
$rng.Font.NameBi = "Arial"

$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Font.NameBi = "Arial"

# ---------------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the end of the "Web/Application
#    Interface" paragraph to wrap the whole "Database Design" body paragraph
#    ("This document includes the data requirements ... or admins.").
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*database design and implementation*") {
        $d.Bookmarks.Add("_GoBack", $p.Range) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Merge the "Hardware" / " Specification" runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Hardware Specification", $true, $false, $false, $false, $false, $true, 1, $false, "Hardware Specification", 2) | Out-Null

Write-Output "done"
